$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 - copy date formatting from A11 (existing date style) then set values
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = 41558
$ws.Range("B12").Value = 2

# Row 13 - copy date formatting from A11 (existing date style) then set values
$ws.Range("A11").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 41559
$ws.Range("B13").Value = 4.5

# Row 28 - sum formula
$ws.Range("B28").Formula = "=SUM(B2:B27)"

# Update selection to mimic end-state cursor position
$ws.Range("B29").Select()
